# Generate Report for Handoff
#
# A new handoff just completed for file "b77c320c-afea-4265-b315-87cf0fd78183",
# so its "Latest Handoff" timestamps are refreshed on every sheet of the
# localization-status report:
#   - Overview!D6        (Latest Handoff Date, row for b77c320c...)
#   - zh-cn!E6            (Latest Handoff Datetime, row for b77c320c...)
#   - de-de!E6            (Latest Handoff Datetime, row for b77c320c...)
#
# These cells are plain text (not real Excel dates), so we assign them as
# strings to keep their original "t=s" (shared-string/text) cell type.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$wsOverview.Range("D6").Value = "2016-32-20 12:32:52"
$wsZhCn.Range("E6").Value     = "2016-03-20 12:32:48"
$wsDeDe.Range("E6").Value     = "2016-03-20 12:32:52"
